$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M21").Value = -36138.332
$ws.Range("K21").Value = 36606.332
$ws.Range("I21").Value = 36606.332
$ws.Range("H21").Value = 36606.332
$ws.Range("K23").Value = 36606.332
$ws.Range("M23").Value = -36372.332
$ws.Range("H23").Value = 36606.332
$ws.Range("I23").Value = 36606.332
$ws.Range("L43").Value = 1933.3334
$ws.Range("J43").Value = 1933.3334
$ws.Range("N43").Value = -2071.3334
$ws.Range("M43").Value = -997.3334
$ws.Range("H43").Value = 1644.3334
$ws.Range("K43").Value = 1066.3334
$ws.Range("I43").Value = 1066.3334
$ws.Range("I138").Value = 5178.1
$ws.Range("K138").Value = 15534.3
$ws.Range("H138").Value = 3026.25
$ws.Range("M138").Value = -10394.3
$ws.Range("L138").Value = 6144.4095
$ws.Range("J138").Value = 2048.1365
$ws.Range("N138").Value = -16424.4095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 3364.5117
$ws.Range("H32").Value = 3878.1228
$ws.Range("I32").Value = 3364.5117
$ws.Range("M32").Value = -3077.5117
$ws.Range("K45").Value = 2028.9
$ws.Range("H45").Value = 2510.7036
$ws.Range("I45").Value = 2028.9
$ws.Range("M45").Value = -1651.9
$ws.Range("K61").Value = 1238.4166
$ws.Range("H61").Value = 2190.7334
$ws.Range("M61").Value = -1026.4166
$ws.Range("I61").Value = 1238.4166
$ws.Range("N109").Value = -70017.71000000001
$ws.Range("H109").Value = 67243.71000000001
$ws.Range("J109").Value = 67243.71000000001
$ws.Range("L109").Value = 67243.71000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -547
$ws.Range("H122").Value = 999
$ws.Range("N122").Value = ""
$ws.Range("J122").Value = 0
$ws.Range("I122").Value = 999
$ws.Range("K122").Value = 2997
$ws.Range("I132").Value = 1376.9445
$ws.Range("K132").Value = 4130.833500000001
$ws.Range("H132").Value = 1739.15
$ws.Range("M132").Value = -1600.833500000001
$ws.Range("M136").Value = -1165.2498
$ws.Range("K136").Value = 3715.2498
$ws.Range("H136").Value = 2190.7334
$ws.Range("I136").Value = 1238.4166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("I22").Value = 0
$ws.Range("M57").Value = ""
$ws.Range("K57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("H80").Value = 9373.817999999999
$ws.Range("N80").Value = -12279.2
$ws.Range("J80").Value = 10283.2
$ws.Range("K80").Value = 280
$ws.Range("L80").Value = 10283.2
$ws.Range("I80").Value = 280
$ws.Range("M80").Value = 718
$ws.Range("K83").Value = 1400
$ws.Range("M83").Value = 3592
$ws.Range("N83").Value = -61400
$ws.Range("H83").Value = 9373.817999999999
$ws.Range("J83").Value = 10283.2
$ws.Range("I83").Value = 280
$ws.Range("L83").Value = 51416
$ws.Range("I86").Value = 1416.6471
$ws.Range("K86").Value = 1416.6471
$ws.Range("M86").Value = -293.6470999999999
$ws.Range("H86").Value = 106741.21
$ws.Range("K89").Value = 7083.2355
$ws.Range("I89").Value = 1416.6471
$ws.Range("M89").Value = -1467.2355
$ws.Range("H89").Value = 106741.21
$ws.Range("M136").Value = ""
$ws.Range("K136").Value = 0
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K16").Value = 953.6667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -666.6667
$ws.Range("H16").Value = 953.6667
$ws.Range("N16").Value = ""
$ws.Range("J16").Value = 0
$ws.Range("I16").Value = 953.6667
$ws.Range("L31").Value = 2069.6924
$ws.Range("J31").Value = 2069.6924
$ws.Range("H31").Value = 1388.1666
$ws.Range("N31").Value = -2659.6924
$ws.Range("N34").Value = -2473.6924
$ws.Range("L34").Value = 2069.6924
$ws.Range("H34").Value = 1388.1666
$ws.Range("J34").Value = 2069.6924
$ws.Range("H62").Value = 7995
$ws.Range("I62").Value = 7995
$ws.Range("K62").Value = 7995
$ws.Range("M62").Value = -7371
$ws.Range("H65").Value = 7995
$ws.Range("M65").Value = -36855
$ws.Range("I65").Value = 7995
$ws.Range("K65").Value = 39975
$ws.Range("M105").Value = -25.71419999999989
$ws.Range("I105").Value = 1772.7142
$ws.Range("H105").Value = 1772.7142
$ws.Range("K105").Value = 1772.7142
$ws.Range("N113").Value = ""
$ws.Range("L113").Value = 0
$ws.Range("H113").Value = 953.6667
$ws.Range("K113").Value = 953.6667
$ws.Range("I113").Value = 953.6667
$ws.Range("M113").Value = 1216.3333
$ws.Range("J113").Value = 0
$ws.Range("I134").Value = 1162.5333
$ws.Range("K134").Value = 3487.5999
$ws.Range("H134").Value = 1744.4
$ws.Range("M134").Value = -952.5999000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J68").Value = 2099.2754
$ws.Range("L68").Value = 6297.8262
$ws.Range("H68").Value = 2047.9028
$ws.Range("N68").Value = -7919.8262
$ws.Range("H71").Value = 2047.9028
$ws.Range("L71").Value = 18893.4786
$ws.Range("N71").Value = -27005.4786
$ws.Range("J71").Value = 2099.2754
$ws.Range("H131").Value = 35766116
$ws.Range("I131").Value = 55555980
$ws.Range("M131").Value = -166662900
$ws.Range("K131").Value = 166667940
$ws.Range("K139").Value = 29316
$ws.Range("H139").Value = 9772
$ws.Range("M139").Value = -24176
$ws.Range("I139").Value = 9772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N20").Value = -11490
$ws.Range("J20").Value = 11000
$ws.Range("H20").Value = 2507333.2
$ws.Range("L20").Value = 11000
$ws.Range("J24").Value = 13200
$ws.Range("N24").Value = -13546
$ws.Range("H24").Value = 2866571.5
$ws.Range("L24").Value = 13200
$ws.Range("N123").Value = -16869.5
$ws.Range("L123").Value = 11969.5
$ws.Range("H123").Value = 11969.5
$ws.Range("J123").Value = 11969.5
$ws.Range("N135").Value = -60270.633
$ws.Range("J135").Value = 50130.633
$ws.Range("H135").Value = 50130.633
$ws.Range("L135").Value = 50130.633

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K16").Value = 6772.5
$ws.Range("L16").Value = 1510.1666
$ws.Range("M16").Value = -6602.5
$ws.Range("H16").Value = 5018.3887
$ws.Range("N16").Value = -1850.1666
$ws.Range("J16").Value = 1510.1666
$ws.Range("I16").Value = 6772.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I132").Value = 679.3
$ws.Range("K132").Value = 2037.9
$ws.Range("H132").Value = 1399.25
$ws.Range("N132").Value = -20057
$ws.Range("L132").Value = 14997
$ws.Range("J132").Value = 4999
$ws.Range("M132").Value = 492.1000000000001
